# Update "想去人数" (interest count) figures for several events across
# the 展览 (Exhibitions), 演出 (Shows) and 全部类型 (All types) sheets.
# Mirrors the site's re-scrape output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions) sheet
$wsExhibit.Range("F6").Value  = 1619
$wsExhibit.Range("F10").Value = 2891
$wsExhibit.Range("F11").Value = 2891
$wsExhibit.Range("F16").Value = 720
$wsExhibit.Range("F18").Value = 6348
$wsExhibit.Range("F19").Value = 250
$wsExhibit.Range("F27").Value = 2489

# 演出 (Shows) sheet
$wsShow.Range("F11").Value = 193

# 全部类型 (All types) sheet
$wsAll.Range("F19").Value = 2891
$wsAll.Range("F23").Value = 193
$wsAll.Range("F26").Value = 720
$wsAll.Range("F27").Value = 6348
$wsAll.Range("F28").Value = 250
$wsAll.Range("F34").Value = 2489

$wb.Save()
